# Update scripts with new TPM values for Slpi-Plscr4 LR-pair sheet.
#
# The original sheet had 6 sender/receiver combinations (rows 2-7):
#   rows 2-4: ECs as sending cluster  -> target clusters ECs, FAPs, MuSCs
#   rows 5-7: MuSCs as sending cluster -> target clusters ECs, FAPs, MuSCs
#
# The updated data drops the "ECs sending cluster" rows entirely and keeps
# only the "MuSCs sending cluster" rows (now renumbered 2-4), refreshed
# with new TPM-derived statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old "ECs" sending-cluster rows (2-4); the old "MuSCs" rows
# (5-7) shift up to become the new rows 2-4.
$ws.Rows("2:4").Delete()

# Row 2: MuSCs -> Slpi -> Plscr4 -> ECs
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.084603
$ws.Range("H2").Value = 0.253809
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("M2").Value = 18.99721533333333
$ws.Range("N2").Value = 56.991646
$ws.Range("O2").Value = 0.3169361933999463
$ws.Range("P2").Value = 0.3169361933999463
$ws.Range("Q2").Value = 1.607221408846
$ws.Range("R2").Value = 14.464992679614
$ws.Range("S2").Value = 0.3169361933999463
$ws.Range("T2").Value = 0.3169361933999463

# Row 3: MuSCs -> Slpi -> Plscr4 -> FAPs
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.084603
$ws.Range("H3").Value = 0.253809
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("N3").Value = 94.79047199999999
$ws.Range("O3").Value = 0.5271392120568722
$ws.Range("P3").Value = 0.5271392120568722
$ws.Range("Q3").Value = 2.673186100872
$ws.Range("R3").Value = 24.058674907848
$ws.Range("S3").Value = 0.5271392120568722
$ws.Range("T3").Value = 0.5271392120568722

# Row 4: MuSCs -> Slpi -> Plscr4 -> MuSCs
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.084603
$ws.Range("H4").Value = 0.253809
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("O4").Value = 0.1559245945431815
$ws.Range("P4").Value = 0.1559245945431815
$ws.Range("Q4").Value = 0.7907123002490001
$ws.Range("R4").Value = 7.116410702241001
$ws.Range("S4").Value = 0.1559245945431815
$ws.Range("T4").Value = 0.1559245945431815
